$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "R package"
$ws.Range("B11").Value = "clusterProfiler"
$ws.Range("C11").Value = "4.10.0"
$ws.Range("D11").Value = "Genetic annotation"

$ws.Range("A12").Value = "Software"
$ws.Range("B12").Value = "SnapGene"
$ws.Range("C12").Value = "8.0.3"
$ws.Range("D12").Value = "Plasmid map viewing and annotation (free version)"

[void]$ws.Range("D13").Select()
